$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cables")
$tbl = $ws.ListObjects.Item(1)

# Rename the "hilos" column header to "hilos usados"
$ws.Range("C1").Value = "hilos usados"

# Add a new table column "Largo original [m]"
$newCol = $tbl.ListColumns.Add()
$ws.Range("D1").Value = "Largo original [m]"

# Match header formatting of the new column to the other header cells
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Correct the "modelo" values (the original file had Conductividad/Oxigeno/PH mismatched)
$ws.Range("B2").Value = "CS511-L"
$ws.Range("B3").Value = "CSIM11-PH-04L"
$ws.Range("B4").Value = "CS547a"

# Fill in the new "Largo original [m]" column data
$ws.Range("D2").Value = 1.4
$ws.Range("D3").Value = 1.4
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 5
$ws.Range("D2:D5").NumberFormat = "0.00"

# Update column widths to fit the new content
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()

# Restore the active selection on the new column
$ws.Range("D2").Select()

$wb.Save()
